$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for "Ají" (Americana (o))
# at Terminal Hortofrutícola Agro Chillán. Insert a new row above the
# existing row 51 so the historic rows (old 51-53) shift down to 52-54,
# then populate the new row 51 with the latest week's data.
$ws.Rows.Item(51).Insert()

$ws.Cells.Item(51, 1).Value  = 7
$ws.Cells.Item(51, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(51, 3).Value  = "Ñuble"
$ws.Cells.Item(51, 4).Value  = 44568
$ws.Cells.Item(51, 5).Value  = 16
$ws.Cells.Item(51, 6).Value  = 100112021
$ws.Cells.Item(51, 7).Value  = "Ají"
$ws.Cells.Item(51, 8).Value  = "Americana (o)"
$ws.Cells.Item(51, 9).Value  = "Primera"
$ws.Cells.Item(51, 10).Value = 100
$ws.Cells.Item(51, 11).Value = 19000
$ws.Cells.Item(51, 12).Value = 20000
$ws.Cells.Item(51, 13).Value = 19500
$ws.Cells.Item(51, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(51, 15).Value = "Región del Maule"
$ws.Cells.Item(51, 16).Value = 1300
$ws.Cells.Item(51, 17).Value = 15
$ws.Cells.Item(51, 18).Value = "Hortaliza"
